$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Grid Analysis")

$ws.Range("C17").Value = 23.9
$ws.Range("F17").Value = 23.9
$ws.Range("G17").Value = 0

# H17/I17 no longer have a detection; set to an explicit empty text value
# (matching the "no detection" pattern used elsewhere in this sheet) rather
# than fully clearing the cells.
$ws.Range("H17").Value = "'"
$ws.Range("I17").Value = "'"
$ws.Range("H17").Style = "Normal"
$ws.Range("I17").Style = "Normal"
